# Refresh cached Universalis market-price snapshot values (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all eight crafting-job sheets, per scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1685.2778
$ws.Range("I58").Value = 302
$ws.Range("J58").Value = 2565.5454
$ws.Range("K58").Value = 906
$ws.Range("L58").Value = 7696.6362
$ws.Range("M58").Value = -756
$ws.Range("N58").Value = -7996.6362

$ws.Range("H64").Value = 3649.475
$ws.Range("J64").Value = 3657.1428
$ws.Range("L64").Value = 3657.1428
$ws.Range("N64").Value = -4153.1428

$ws.Range("H67").Value = 3649.475
$ws.Range("J67").Value = 3657.1428
$ws.Range("L67").Value = 3657.1428
$ws.Range("N67").Value = -5373.1428

$ws.Range("H76").Value = 6733.5835
$ws.Range("I76").Value = 7475.375
$ws.Range("J76").Value = 5250
$ws.Range("K76").Value = 7475.375
$ws.Range("L76").Value = 5250
$ws.Range("M76").Value = -7160.375
$ws.Range("N76").Value = -5880

$ws.Range("H79").Value = 6733.5835
$ws.Range("I79").Value = 7475.375
$ws.Range("J79").Value = 5250
$ws.Range("K79").Value = 7475.375
$ws.Range("L79").Value = 5250
$ws.Range("M79").Value = -6383.375
$ws.Range("N79").Value = -7434

$ws.Range("H82").Value = 10601.462
$ws.Range("I82").Value = 603.8
$ws.Range("J82").Value = 16850
$ws.Range("K82").Value = 1811.4
$ws.Range("L82").Value = 50550
$ws.Range("M82").Value = -1405.4
$ws.Range("N82").Value = -51362

$ws.Range("H85").Value = 10601.462
$ws.Range("I85").Value = 603.8
$ws.Range("J85").Value = 16850
$ws.Range("K85").Value = 1811.4
$ws.Range("L85").Value = 50550
$ws.Range("M85").Value = -407.3999999999999
$ws.Range("N85").Value = -53358

$ws.Range("H93").Value = 178571.42
$ws.Range("J93").Value = 178571.42
$ws.Range("L93").Value = 178571.42
$ws.Range("N93").Value = -183563.42

$ws.Range("H109").Value = 54117.875
$ws.Range("J109").Value = 54117.875
$ws.Range("L109").Value = 54117.875
$ws.Range("N109").Value = -56891.875

$ws.Range("H115").Value = 2666.25
$ws.Range("I115").Value = 1895
$ws.Range("J115").Value = 4980
$ws.Range("K115").Value = 5685
$ws.Range("L115").Value = 14940
$ws.Range("M115").Value = -4118
$ws.Range("N115").Value = -18074

$ws.Range("H129").Value = 1045.024
$ws.Range("J129").Value = 1081.6075
$ws.Range("L129").Value = 3244.8225
$ws.Range("N129").Value = -13244.8225

$ws.Range("H135").Value = 565.2857
$ws.Range("I135").Value = 401.16
$ws.Range("K135").Value = 3610.44
$ws.Range("M135").Value = -1075.44

$ws.Range("H137").Value = 2265.2
$ws.Range("I137").Value = 1746.619
$ws.Range("J137").Value = 3043.0715
$ws.Range("K137").Value = 5239.857
$ws.Range("L137").Value = 9129.2145
$ws.Range("M137").Value = -2689.857
$ws.Range("N137").Value = -14229.2145

$ws.Range("H138").Value = 3490.0544
$ws.Range("I138").Value = 1188.0344
$ws.Range("J138").Value = 6057.6924
$ws.Range("K138").Value = 3564.1032
$ws.Range("L138").Value = 18173.0772
$ws.Range("M138").Value = 1575.8968
$ws.Range("N138").Value = -28453.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28472.73
$ws.Range("I32").Value = 39926.53
$ws.Range("K32").Value = 39926.53
$ws.Range("M32").Value = -39639.53

$ws.Range("H61").Value = 1887.9348
$ws.Range("I61").Value = 1763.8605
$ws.Range("K61").Value = 1763.8605
$ws.Range("M61").Value = -1551.8605

$ws.Range("H74").Value = 1137.2538
$ws.Range("I74").Value = 1008.6458
$ws.Range("J74").Value = 1462.1578
$ws.Range("K74").Value = 1008.6458
$ws.Range("L74").Value = 1462.1578
$ws.Range("M74").Value = -134.6458
$ws.Range("N74").Value = -3210.1578

$ws.Range("H77").Value = 1137.2538
$ws.Range("I77").Value = 1008.6458
$ws.Range("J77").Value = 1462.1578
$ws.Range("K77").Value = 5043.229
$ws.Range("L77").Value = 7310.789
$ws.Range("M77").Value = -675.2290000000003
$ws.Range("N77").Value = -16046.789

$ws.Range("H132").Value = 4265.1113
$ws.Range("I132").Value = 4164.6904
$ws.Range("J132").Value = 5671
$ws.Range("K132").Value = 12494.0712
$ws.Range("L132").Value = 17013
$ws.Range("M132").Value = -9964.071200000002
$ws.Range("N132").Value = -22073

$ws.Range("H136").Value = 1887.9348
$ws.Range("I136").Value = 1763.8605
$ws.Range("K136").Value = 5291.5815
$ws.Range("M136").Value = -2741.5815

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2832.6155
$ws.Range("I134").Value = 3024.889
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 9074.667000000001
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -6539.667000000001
$ws.Range("N134").Value = -12270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1655.5753
$ws.Range("I31").Value = 2138
$ws.Range("J31").Value = 1388.7021
$ws.Range("K31").Value = 2138
$ws.Range("L31").Value = 1388.7021
$ws.Range("M31").Value = -1843
$ws.Range("N31").Value = -1978.7021

$ws.Range("H34").Value = 1655.5753
$ws.Range("I34").Value = 2138
$ws.Range("J34").Value = 1388.7021
$ws.Range("K34").Value = 2138
$ws.Range("L34").Value = 1388.7021
$ws.Range("M34").Value = -1936
$ws.Range("N34").Value = -1792.7021

$ws.Range("H132").Value = 3072.162
$ws.Range("I132").Value = 2666.8386
$ws.Range("J132").Value = 5166.3335
$ws.Range("K132").Value = 8000.5158
$ws.Range("L132").Value = 15499.0005
$ws.Range("M132").Value = -5470.5158
$ws.Range("N132").Value = -20559.0005

$ws.Range("H134").Value = 3332.6667
$ws.Range("I134").Value = 2260
$ws.Range("K134").Value = 6780
$ws.Range("M134").Value = -4245

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 740.2969000000001
$ws.Range("I68").Value = 550.46155
$ws.Range("K68").Value = 1651.38465
$ws.Range("M68").Value = -840.38465

$ws.Range("H71").Value = 740.2969000000001
$ws.Range("I71").Value = 550.46155
$ws.Range("K71").Value = 4954.15395
$ws.Range("M71").Value = -898.1539499999999

$ws.Range("H126").Value = 5600
$ws.Range("I126").Value = 4240
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 12720
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -7780
$ws.Range("N126").Value = -36880

$ws.Range("H139").Value = 2401.625
$ws.Range("I139").Value = 2181.5386
$ws.Range("J139").Value = 3355.3333
$ws.Range("K139").Value = 6544.6158
$ws.Range("L139").Value = 10065.9999
$ws.Range("M139").Value = -1404.6158
$ws.Range("N139").Value = -20345.9999

$ws.Range("H140").Value = 8158.3335
$ws.Range("I140").Value = 1030
$ws.Range("J140").Value = 9049.375
$ws.Range("K140").Value = 3090
$ws.Range("L140").Value = 27148.125
$ws.Range("M140").Value = 2090
$ws.Range("N140").Value = -37508.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2215.879
$ws.Range("I80").Value = 2240.238
$ws.Range("J80").Value = 2173.25
$ws.Range("K80").Value = 2240.238
$ws.Range("L80").Value = 2173.25
$ws.Range("M80").Value = -1242.238
$ws.Range("N80").Value = -4169.25

$ws.Range("H83").Value = 2215.879
$ws.Range("I83").Value = 2240.238
$ws.Range("J83").Value = 2173.25
$ws.Range("K83").Value = 11201.19
$ws.Range("L83").Value = 10866.25
$ws.Range("M83").Value = -6209.189999999999
$ws.Range("N83").Value = -20850.25

$ws.Range("H132").Value = 2350.325
$ws.Range("I132").Value = 1869.4828
$ws.Range("K132").Value = 5608.4484
$ws.Range("M132").Value = -3078.4484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4649.1177
$ws.Range("I132").Value = 4775.92
$ws.Range("J132").Value = 4296.8887
$ws.Range("K132").Value = 14327.76
$ws.Range("L132").Value = 12890.6661
$ws.Range("M132").Value = -11797.76
$ws.Range("N132").Value = -17950.6661

$ws.Range("H136").Value = 1513.5
$ws.Range("I136").Value = 895.25
$ws.Range("K136").Value = 2685.75
$ws.Range("M136").Value = -135.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 7966.6665
$ws.Range("J69").Value = 7966.6665
$ws.Range("L69").Value = 7966.6665
$ws.Range("N69").Value = -9464.666499999999

$ws.Range("H72").Value = 7966.6665
$ws.Range("J72").Value = 7966.6665
$ws.Range("L72").Value = 23899.9995
$ws.Range("N72").Value = -31387.9995

$ws.Range("H132").Value = 2594.3726
$ws.Range("I132").Value = 2113.5264
$ws.Range("J132").Value = 3999.923
$ws.Range("K132").Value = 6340.5792
$ws.Range("L132").Value = 11999.769
$ws.Range("M132").Value = -3810.5792
$ws.Range("N132").Value = -17059.769
